# Commit: "Make elimination language independent" (Close #406)
#
# On the "Variables" sheet, the elimination value for each variable used to be
# stored as a translated label in the "fo_elimination" column (e.g. "I alt",
# "Fólkatalið fyrst í árinum"). It is now stored as a language-independent
# CODE in a new "elimination" column, which is positioned right after
# "timeval" (i.e. before "fo_variable-label" rather than after it).
#
# Net effect on the "Variables" table column layout:
#   old: pivot, order, variable-code, variable-type, timeval,
#        fo_variable-label, fo_elimination, fo_note
#   new: pivot, order, variable-code, variable-type, timeval,
#        elimination, fo_variable-label, fo_note
#
# i.e. column F ("fo_variable-label") shifts right into column G, column G
# ("fo_elimination") is replaced by the new language-independent "elimination"
# column values, and column H ("fo_note") is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# --- capture the current (pre-edit) column F values before overwriting ---
$fVarLabelHeader = $ws.Range("F1").Value2
$fVarLabel2 = $ws.Range("F2").Value2
$fVarLabel3 = $ws.Range("F3").Value2
$fVarLabel4 = $ws.Range("F4").Value2
$fVarLabel5 = $ws.Range("F5").Value2

# --- shift the old "fo_variable-label" column (F) into column G ---
$ws.Range("G1").Value = $fVarLabelHeader
$ws.Range("G2").Value = $fVarLabel2
$ws.Range("G3").Value = $fVarLabel3
$ws.Range("G4").Value = $fVarLabel4
$ws.Range("G5").Value = $fVarLabel5

# --- write the new language-independent "elimination" column into F ---
$ws.Range("F1").Value = "elimination"
$ws.Range("F2").Value = "T"
$ws.Range("F3").ClearContents()
$ws.Range("F4").Value = "P"
$ws.Range("F5").ClearContents()

# --- update selections to match the saved state ---
$wsCells = $wb.Worksheets.Item("Cells")
$wsCells.Activate()
$wsCells.Range("B2").Select()

$ws.Activate()
$ws.Range("F4").Select()
